$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").ClearContents()
$ws.Range("F11").Select()
